$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 1500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1385
$ws.Range("N5").ClearContents()

$ws.Range("H46").Value = 1063.3334
$ws.Range("J46").Value = 595
$ws.Range("L46").Value = 1785
$ws.Range("N46").Value = -2023

$ws.Range("H60").Value = 1063.3334
$ws.Range("J60").Value = 595
$ws.Range("L60").Value = 1785
$ws.Range("N60").Value = -2753

$ws.Range("H92").Value = 423.94116
$ws.Range("I92").Value = 351.92856
$ws.Range("K92").Value = 351.92856
$ws.Range("M92").Value = 896.0714399999999

$ws.Range("H107").Value = 406.4737
$ws.Range("J107").Value = 777
$ws.Range("L107").Value = 777
$ws.Range("N107").Value = -4617

$ws.Range("H114").Value = 61500
$ws.Range("J114").Value = 61500
$ws.Range("L114").Value = 61500
$ws.Range("N114").Value = -70178

$ws.Range("H125").Value = 100
$ws.Range("I125").Value = 100
$ws.Range("K125").Value = 900
$ws.Range("M125").Value = 1560

$ws.Range("H137").Value = 1419.2
$ws.Range("J137").Value = 2874
$ws.Range("L137").Value = 8622
$ws.Range("N137").Value = -13722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2527.5715
$ws.Range("J45").Value = 1897
$ws.Range("L45").Value = 1897
$ws.Range("N45").Value = -2651

$ws.Range("H61").Value = 3999.5
$ws.Range("I61").Value = 3999.5
$ws.Range("K61").Value = 3999.5
$ws.Range("M61").Value = -3787.5

$ws.Range("H110").Value = 999.5
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 999
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 999
$ws.Range("M110").Value = 1045
$ws.Range("N110").Value = -5089

$ws.Range("H132").Value = 2959.8
$ws.Range("I132").Value = 2849.75
$ws.Range("K132").Value = 8549.25
$ws.Range("M132").Value = -6019.25

$ws.Range("H136").Value = 3999.5
$ws.Range("I136").Value = 3999.5
$ws.Range("K136").Value = 11998.5
$ws.Range("M136").Value = -9448.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 3688.25
$ws.Range("I75").Value = 3688.25
$ws.Range("K75").Value = 3688.25
$ws.Range("M75").Value = -2752.25

$ws.Range("H78").Value = 3688.25
$ws.Range("I78").Value = 3688.25
$ws.Range("K78").Value = 11064.75
$ws.Range("M78").Value = -6384.75

$ws.Range("H86").Value = 2332
$ws.Range("I86").Value = 2198.6
$ws.Range("K86").Value = 2198.6
$ws.Range("M86").Value = -1075.6

$ws.Range("H89").Value = 2332
$ws.Range("I89").Value = 2198.6
$ws.Range("K89").Value = 10993
$ws.Range("M89").Value = -5377

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900.5
$ws.Range("I16").Value = 900.5
$ws.Range("K16").Value = 900.5
$ws.Range("M16").Value = -613.5

$ws.Range("H31").Value = 4788.6924
$ws.Range("I31").Value = 2796.375
$ws.Range("K31").Value = 2796.375
$ws.Range("M31").Value = -2501.375

$ws.Range("H34").Value = 4788.6924
$ws.Range("I34").Value = 2796.375
$ws.Range("K34").Value = 2796.375
$ws.Range("M34").Value = -2594.375

$ws.Range("H62").Value = 3800
$ws.Range("I62").Value = 3750
$ws.Range("K62").Value = 3750
$ws.Range("M62").Value = -3126

$ws.Range("H65").Value = 3800
$ws.Range("I65").Value = 3750
$ws.Range("K65").Value = 18750
$ws.Range("M65").Value = -15630

$ws.Range("H113").Value = 900.5
$ws.Range("I113").Value = 900.5
$ws.Range("K113").Value = 900.5
$ws.Range("M113").Value = 1269.5

$ws.Range("H134").Value = 1999.5
$ws.Range("I134").Value = 1999
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5997
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -3462
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 522.375
$ws.Range("I33").Value = 617.8
$ws.Range("J33").Value = 363.33334
$ws.Range("K33").Value = 3706.8
$ws.Range("L33").Value = 2180.00004
$ws.Range("M33").Value = -3423.8
$ws.Range("N33").Value = -2746.00004

$ws.Range("H117").Value = 66.333336
$ws.Range("I117").Value = 66.333336
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 199.000008
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 3242.999992
$ws.Range("N117").ClearContents()

$ws.Range("H137").Value = 988
$ws.Range("I137").Value = 988
$ws.Range("K137").Value = 2964
$ws.Range("M137").Value = 2136

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 42500
$ws.Range("J15").Value = 42500
$ws.Range("L15").Value = 42500
$ws.Range("N15").Value = -43076

$ws.Range("H43").Value = 11629.167
$ws.Range("I43").Value = 3943.75
$ws.Range("J43").Value = 27000
$ws.Range("K43").Value = 3943.75
$ws.Range("L43").Value = 27000
$ws.Range("M43").Value = -3792.75
$ws.Range("N43").Value = -27302

$ws.Range("H81").Value = 42500
$ws.Range("J81").Value = 42500
$ws.Range("L81").Value = 42500
$ws.Range("N81").Value = -44496

$ws.Range("H84").Value = 42500
$ws.Range("J84").Value = 42500
$ws.Range("L84").Value = 127500
$ws.Range("N84").Value = -137484

$ws.Range("H122").Value = 5664.9414
$ws.Range("J122").Value = 8664.833000000001
$ws.Range("L122").Value = 25994.499
$ws.Range("N122").Value = -30894.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8838.5
$ws.Range("I7").Value = 9224.5
$ws.Range("J7").Value = 8452.5
$ws.Range("K7").Value = 9224.5
$ws.Range("L7").Value = 8452.5
$ws.Range("M7").Value = -9112.5
$ws.Range("N7").Value = -8676.5

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 1750
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876

$ws.Range("H126").Value = 8838.5
$ws.Range("I126").Value = 9224.5
$ws.Range("J126").Value = 8452.5
$ws.Range("K126").Value = 27673.5
$ws.Range("L126").Value = 25357.5
$ws.Range("M126").Value = -25203.5
$ws.Range("N126").Value = -30297.5

